$d = $word.ActiveDocument

# Locate the "Desired Position:" label so the rest of the edits can be scoped
# precisely to its neighbouring value cell, instead of matching look-alike
# text located elsewhere in the CV (e.g. "Full Stack Software Engineer").
$label = $d.Content
$label.Find.Execute("Desired Position:", $false)
$labelEnd = $label.End

# 1) "Desired Position" value: "Senior Software Engineer / Team Lead" -> "Frond End / Full Stack"
$target = $d.Range($labelEnd, $labelEnd + 60)
$target.Find.Execute("Senior Software Engineer / Team Lead", $true, $false, $false, $false, $false,
                      $true, 1, $false, "Frond End / Full Stack", 2)

# Recreate the hidden "_GoBack" bookmark at the point of this latest edit
# (between "Full S" and "tack"), mirroring Word's own behaviour of moving
# that bookmark to the location of the most recent edit.
$scope = $d.Range($labelEnd, $labelEnd + 40)
$scope.Find.Execute("Full S", $false)
$editPos = $scope.End
$bmRange = $d.Range($editPos, $editPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 2) "Integration and develo|pment of ERP system UNA.md (http://una.md/en)" was split into
#    two runs around the old "_GoBack" bookmark. Re-applying the full, correctly spelled
#    phrase merges the runs back together and (because the edit touches that spot) drops
#    the stale bookmark that used to live there.
$d.Content.Find.Execute("Integration and development of ERP system UNA.md (http://una.md/en)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Integration and development of ERP system UNA.md (http://una.md/en)", 2)
